{"js": "// The document contains a single 20-row x 5-column table of simple\n// arithmetic \"question=answer\" strings (one per cell, one run per\n// paragraph). The edit replaces every cell's text with a new\n// expression while keeping the table's row/column layout and the\n// existing run/paragraph formatting untouched.\n\nconst newValues = [\n  [\"94-77=17\", \"9+56=65\", \"19+68=87\", \"54-5=49\", \"63-37=26\"],\n  [\"31-28=3\", \"39+33=72\", \"45-19=26\", \"70-55=15\", \"87+8=95\"],\n  [\"27+24=51\", \"29+36=65\", \"52-49=3\", \"50-33=17\", \"10-7=3\"],\n  [\"28+49=77\", \"9+29=38\", \"49+43=92\", \"52-24=28\", \"90-9=81\"],\n  [\"19+69=88\", \"9+47=56\", \"80-75=5\", \"48+47=95\", \"65-36=29\"],\n  [\"37+6=43\", \"50-45=5\", \"60-18=42\", \"82-39=43\", \"87-58=29\"],\n  [\"50-14=36\", \"42-15=27\", \"90-81=9\", \"27+28=55\", \"41-6=35\"],\n  [\"84-27=57\", \"70-27=43\", \"62-35=27\", \"67+6=73\", \"84-8=76\"],\n  [\"23-5=18\", \"59+6=65\", \"81-13=68\", \"34+7=41\", \"58+34=92\"],\n  [\"90-52=38\", \"14+39=53\", \"32-18=14\", \"28+3=31\", \"30-11=19\"],\n  [\"61-52=9\", \"83-36=47\", \"28+37=65\", \"94-75=19\", \"71-59=12\"],\n  [\"7+77=84\", \"74-55=19\", \"90-2=88\", \"96-28=68\", \"49+19=68\"],\n  [\"57+38=95\", \"35+49=84\", \"56+26=82\", \"16+17=33\", \"4+79=83\"],\n  [\"37-8=29\", \"3+88=91\", \"72-48=24\", \"94-19=75\", \"19+53=72\"],\n  [\"4+88=92\", \"40-18=22\", \"80-2=78\", \"26+35=61\", \"16+67=83\"],\n  [\"42+19=61\", \"90-83=7\", \"26+6=32\", \"7+9=16\", \"18+58=76\"],\n  [\"70-66=4\", \"93-67=26\", \"5+77=82\", \"53-34=19\", \"58+38=96\"],\n  [\"9+57=66\", \"91-38=53\", \"6+77=83\", \"25+46=71\", \"28+29=57\"],\n  [\"28+38=66\", \"52+9=61\", \"92-45=47\", \"86-48=38\", \"15+78=93\"],\n  [\"26-9=17\", \"9+42=51\", \"83-29=54\", \"66-19=47\", \"79+9=88\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Update cell-by-cell via TableCell.value so only the text content\n// changes; existing paragraph/run formatting (font, size, alignment)\n// already present in each cell is left exactly as-is.\nfor (let r = 0; r < newValues.length && r < table.rowCount; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = newValues[r][c];\n  }\n}\nawait context.sync();\n", "ps1": "# The document contains a single 20-row x 5-column table of simple\n# arithmetic \"question=answer\" strings (one per cell). The edit replaces\n# every cell's text with a new expression while keeping the table's\n# row/column layout and each cell's existing paragraph/run formatting\n# untouched.\n\n$newValues = @(\n    @(\"94-77=17\", \"9+56=65\", \"19+68=87\", \"54-5=49\", \"63-37=26\"),\n    @(\"31-28=3\", \"39+33=72\", \"45-19=26\", \"70-55=15\", \"87+8=95\"),\n    @(\"27+24=51\", \"29+36=65\", \"52-49=3\", \"50-33=17\", \"10-7=3\"),\n    @(\"28+49=77\", \"9+29=38\", \"49+43=92\", \"52-24=28\", \"90-9=81\"),\n    @(\"19+69=88\", \"9+47=56\", \"80-75=5\", \"48+47=95\", \"65-36=29\"),\n    @(\"37+6=43\", \"50-45=5\", \"60-18=42\", \"82-39=43\", \"87-58=29\"),\n    @(\"50-14=36\", \"42-15=27\", \"90-81=9\", \"27+28=55\", \"41-6=35\"),\n    @(\"84-27=57\", \"70-27=43\", \"62-35=27\", \"67+6=73\", \"84-8=76\"),\n    @(\"23-5=18\", \"59+6=65\", \"81-13=68\", \"34+7=41\", \"58+34=92\"),\n    @(\"90-52=38\", \"14+39=53\", \"32-18=14\", \"28+3=31\", \"30-11=19\"),\n    @(\"61-52=9\", \"83-36=47\", \"28+37=65\", \"94-75=19\", \"71-59=12\"),\n    @(\"7+77=84\", \"74-55=19\", \"90-2=88\", \"96-28=68\", \"49+19=68\"),\n    @(\"57+38=95\", \"35+49=84\", \"56+26=82\", \"16+17=33\", \"4+79=83\"),\n    @(\"37-8=29\", \"3+88=91\", \"72-48=24\", \"94-19=75\", \"19+53=72\"),\n    @(\"4+88=92\", \"40-18=22\", \"80-2=78\", \"26+35=61\", \"16+67=83\"),\n    @(\"42+19=61\", \"90-83=7\", \"26+6=32\", \"7+9=16\", \"18+58=76\"),\n    @(\"70-66=4\", \"93-67=26\", \"5+77=82\", \"53-34=19\", \"58+38=96\"),\n    @(\"9+57=66\", \"91-38=53\", \"6+77=83\", \"25+46=71\", \"28+29=57\"),\n    @(\"28+38=66\", \"52+9=61\", \"92-45=47\", \"86-48=38\", \"15+78=93\"),\n    @(\"26-9=17\", \"9+42=51\", \"83-29=54\", \"66-19=47\", \"79+9=88\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
